$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.287196755409241
$ws.Range("B1").Value = 2.152289390563965
$ws.Range("C1").Value = 4.790530204772949
$ws.Range("D1").Value = 3.275928497314453
$ws.Range("E1").Value = 1.358094811439514
